$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.005.73'
$ws.Range("E2").Value = '  +1.89%  '

# Row 3
$ws.Range("D3").Value = '2.526.38'
$ws.Range("E3").Value = '  +0.21%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.76'
$ws.Range("E5").Value = '  +1.67%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.49'
$ws.Range("E6").Value = '  +2.72%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.531'
$ws.Range("E8").Value = '  +0.63%  '

# Row 9
$ws.Range("D9").Value = '2.525.22'
$ws.Range("E9").Value = '  +0.21%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.141'
$ws.Range("E10").Value = '  +0.68%  '

# Row 11
$ws.Range("E11").Value = '  +2.66%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.11'
$ws.Range("E12").Value = '  -0.55%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.345'
$ws.Range("E13").Value = '  -1.41%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.71'
$ws.Range("E14").Value = '  -0.28%  '

# Row 15
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.990.23'
$ws.Range("E15").Value = '  +0.72%  '

# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000179'
$ws.Range("E16").Value = '  +1.19%  '

# Row 17
$ws.Range("D17").Value = '67.908.38'
$ws.Range("E17").Value = '  +2.04%  '

# Row 18
$ws.Range("D18").Value = '2.521.96'
$ws.Range("E18").Value = '  +0.13%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.07'
$ws.Range("E19").Value = '  +2.13%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.67'
$ws.Range("E20").Value = '  +2.74%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '366.11'
$ws.Range("E21").Value = '  +4.90%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.20'
$ws.Range("E22").Value = '  -0.44%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.68'
$ws.Range("E23").Value = '  +0.45%  '

# Row 24
$ws.Range("E24").Value = '  -0.07%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.93'
$ws.Range("E25").Value = '  -2.12%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.92'
$ws.Range("E26").Value = '  +1.33%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.14'
$ws.Range("E27").Value = '  +3.03%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.17%  '

# Row 29
$ws.Range("D29").Value = '2.657.37'
$ws.Range("E29").Value = '  +0.49%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0993'
$ws.Range("E30").Value = '  +0.65%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.35'
$ws.Range("E31").Value = '  +1.08%  '

# Row 32
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '534.95'
$ws.Range("E32").Value = '  +1.25%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.33'
$ws.Range("E33").Value = '  +0.48%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.88'
$ws.Range("E34").Value = '  +2.37%  '

# Row 35
$ws.Range("E35").Value = '  -1.25%  '

# Row 36
$ws.Range("E36").Value = '  +0.04%  '

# Row 37
$ws.Range("E37").Value = '  -0.23%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '157.45'
$ws.Range("E38").Value = '  +0.15%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.80'
$ws.Range("E39").Value = '  +0.74%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.68'
$ws.Range("E40").Value = '  +1.52%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.354'
$ws.Range("E41").Value = '  -0.43%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.19'
$ws.Range("E42").Value = '  +1.84%  '

# Row 43
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.80'
$ws.Range("E43").Value = '  -0.19%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.54'
$ws.Range("E44").Value = '  +0.59%  '

# Row 45
$ws.Range("E45").Value = '  -0.08%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '147.28'
$ws.Range("E46").Value = '  -1.93%  '

# Row 47
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₆0280'
$ws.Range("E47").Value = '  +4.18%  '

# Row 48
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.559'
$ws.Range("E48").Value = '  -0.08%  '

# Row 49
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.73'
$ws.Range("E49").Value = '  +0.83%  '

# Row 50
$ws.Range("E50").Value = '  -1.58%  '

# Row 51
$ws.Range("E51").Value = '  +0.24%  '
